$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updating odds values per the 2026-01-02 Betfair data refresh
# Row 2
$ws.Range("G2").Value = 1.84
$ws.Range("P2").Value = 2.22
$ws.Range("T2").Value = 1.72
$ws.Range("W2").Value = 2.18
$ws.Range("AC2").Value = 9
$ws.Range("AG2").Value = 9.800000000000001

# Row 3
$ws.Range("F3").Value = 3.8
$ws.Range("G3").Value = 4.5
$ws.Range("H3").Value = 1.91
$ws.Range("I3").Value = 2.12
$ws.Range("J3").Value = 3.7
$ws.Range("N3").Value = 4.4
$ws.Range("P3").Value = 2.16
$ws.Range("Q3").Value = 1.7
$ws.Range("R3").Value = 1.47
$ws.Range("S3").Value = 2.72
$ws.Range("T3").Value = 1.64
$ws.Range("U3").Value = 2.24
$ws.Range("V3").Value = 1.89
$ws.Range("W3").Value = 1.29
$ws.Range("X3").Value = 20
$ws.Range("Y3").Value = 12.5
$ws.Range("AA3").Value = 24
$ws.Range("AB3").Value = 18.5
$ws.Range("AC3").Value = 9.4
$ws.Range("AH3").Value = 17
$ws.Range("AO3").Value = 11.5

# Row 4
$ws.Range("F4").Value = 1.51
$ws.Range("G4").Value = 1.63
$ws.Range("H4").Value = 6.6
$ws.Range("I4").Value = 8
$ws.Range("J4").Value = 4.1
$ws.Range("K4").Value = 4.9
$ws.Range("N4").Value = 3.7
$ws.Range("O4").Value = 1.29
$ws.Range("P4").Value = 1.93
$ws.Range("Q4").Value = 1.89
$ws.Range("R4").Value = 1.34
$ws.Range("S4").Value = 3.3
$ws.Range("T4").Value = 1.98
$ws.Range("U4").Value = 1.85
$ws.Range("V4").Value = 1.14
$ws.Range("W4").Value = 2.6
$ws.Range("X4").Value = 990
$ws.Range("AA4").Value = 270
$ws.Range("AB4").Value = 9.4
$ws.Range("AC4").Value = 10.5
$ws.Range("AJ4").Value = 15
$ws.Range("AK4").Value = 1000
$ws.Range("AM4").Value = 180
$ws.Range("AN4").Value = 9.6

# Row 5
$ws.Range("G5").Value = 290
$ws.Range("H5").Value = 1.13
$ws.Range("I5").Value = 12
$ws.Range("J5").Value = 1.09
$ws.Range("V5").Value = 1.09

# Row 7
$ws.Range("F7").Value = 4.9
$ws.Range("G7").Value = 5.8
$ws.Range("H7").Value = 1.58
$ws.Range("I7").Value = 1.63
$ws.Range("J7").Value = 4.7
$ws.Range("K7").Value = 5.4
$ws.Range("N7").Value = 6.2
$ws.Range("O7").Value = 1.15
$ws.Range("P7").Value = 2.78
$ws.Range("Q7").Value = 1.39
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 2.1
$ws.Range("T7").Value = 1.56
$ws.Range("U7").Value = 2.42
$ws.Range("V7").Value = 2.56
$ws.Range("W7").Value = 1.21
$ws.Range("X7").Value = 990
$ws.Range("Z7").Value = 14
$ws.Range("AA7").Value = 20
$ws.Range("AE7").Value = 16
$ws.Range("AG7").Value = 980
$ws.Range("AH7").Value = 23
$ws.Range("AJ7").Value = 140
$ws.Range("AN7").Value = 1000
$ws.Range("AO7").Value = 5.9

# Row 8
$ws.Range("F8").Value = 8
$ws.Range("G8").Value = 8.6
$ws.Range("H8").Value = 1.5
$ws.Range("M8").Value = 1.06
$ws.Range("O8").Value = 1.33
$ws.Range("S8").Value = 3.55
$ws.Range("T8").Value = 2.14
$ws.Range("U8").Value = 1.78
$ws.Range("W8").Value = 1.13
$ws.Range("AL8").Value = 140
$ws.Range("AN8").Value = 220

# Row 9
$ws.Range("H9").Value = 6.4

# Row 10
$ws.Range("G10").Value = 2.94
$ws.Range("H10").Value = 2.68
$ws.Range("N10").Value = 3.8

# Row 11
$ws.Range("H11").Value = 1.58
$ws.Range("S11").Value = 3.55
$ws.Range("AC11").Value = 9.199999999999999

# Row 12
$ws.Range("F12").Value = 2.18
$ws.Range("G12").Value = 2.22
$ws.Range("H12").Value = 4.6
$ws.Range("I12").Value = 4.7
$ws.Range("Q12").Value = 3.3
$ws.Range("V12").Value = 1.27
$ws.Range("W12").Value = 1.81
$ws.Range("AF12").Value = 10.5
$ws.Range("AM12").Value = 350
$ws.Range("AN12").Value = 40
$ws.Range("AO12").Value = 200

# Row 13
$ws.Range("F13").Value = 1.88
$ws.Range("G13").Value = 1.9
$ws.Range("H13").Value = 4.9
$ws.Range("I13").Value = 5.3
$ws.Range("N13").Value = 3.25
$ws.Range("Q13").Value = 2.12
$ws.Range("R13").Value = 1.29
$ws.Range("S13").Value = 3.95
$ws.Range("T13").Value = 1.96
$ws.Range("V13").Value = 1.24
$ws.Range("W13").Value = 2.1
$ws.Range("X13").Value = 12
$ws.Range("Y13").Value = 15.5
$ws.Range("AA13").Value = 130
$ws.Range("AB13").Value = 7.8
$ws.Range("AC13").Value = 8.199999999999999
$ws.Range("AD13").Value = 20
$ws.Range("AE13").Value = 75
$ws.Range("AG13").Value = 10.5
$ws.Range("AI13").Value = 85
$ws.Range("AJ13").Value = 21
$ws.Range("AL13").Value = 44
$ws.Range("AM13").Value = 140
